# convection.xlsx — "Add files via upload" commit
#
# Net effect (reconstructed from the OOXML diff):
#   * Sheet "2_" (was empty) gets a new 5-row multiple-choice quiz block
#     about how surface temperature / overall flux change as convection
#     becomes more effective.
#   * Sheet "3_" gets its old quiz block (about the discretized 2nd
#     derivative) replaced with a new 5-row quiz block about how the
#     temperature gradients *near the wall* change as convection becomes
#     more effective.
#   * The previously-selected/active sheet "1_" loses its tab selection
#     and its in-sheet selection collapses to A1:C5 (no sticky ActiveCell).
#   * Sheet "3_" becomes the active tab/sheet, with C3 selected.
#
# (Sheets "0_" and "1_" keep exactly the same visible text — only their
#  underlying shared-string ids shift because six now-unused strings were
#  removed from the shared string table elsewhere, so nothing to touch
#  there beyond the selection bookkeeping above.)

$wb = $excel.ActiveWorkbook

$ws0 = $wb.Worksheets.Item("0_")
$ws1 = $wb.Worksheets.Item("1_")
$ws2 = $wb.Worksheets.Item("2_")
$ws3 = $wb.Worksheets.Item("3_")

# ---------------------------------------------------------------------
# Sheet "2_" : new quiz content (previously completely empty)
# ---------------------------------------------------------------------
$ws2.Cells.Item(1,1).Value = "According to the plots, what happens to the temperature of the surface and overall flux as the convection becomes more effective?"
$ws2.Cells.Item(1,2).Value = "Correct"
$ws2.Cells.Item(1,3).Value = "Comment"
$ws2.Rows.Item(1).RowHeight = 90

$ws2.Cells.Item(2,1).Value = "The surface temperature increases, and the flux increases"
$ws2.Cells.Item(2,2).Value = "N"
$ws2.Rows.Item(2).RowHeight = 45

$ws2.Cells.Item(3,1).Value = "The surface temperature increases, but the flux goes down"
$ws2.Cells.Item(3,2).Value = "N"
$ws2.Rows.Item(3).RowHeight = 45

$ws2.Cells.Item(4,1).Value = "The surface temperature decreases, but the flux goes up"
$ws2.Cells.Item(4,2).Value = "Y"
$ws2.Cells.Item(4,3).Value = "Yep!  As the convection becomes more violent, it throws colder fluid next to the surface, cooling the surface.  The cooler surface temperature increases the gradient in the solid, which leads to a larger flux of energy out of the solid."
$ws2.Rows.Item(4).RowHeight = 120

$ws2.Cells.Item(5,1).Value = "The surface temperature decreases, and the flux goes down"
$ws2.Cells.Item(5,2).Value = "N"
$ws2.Rows.Item(5).RowHeight = 45

# ---------------------------------------------------------------------
# Sheet "3_" : replace old quiz content with the new gradient-near-wall quiz
# ---------------------------------------------------------------------
$ws3.Cells.Item(1,1).Value = "According to the plots, what happens to the temperature gradients *near the wall* as convection becomes more effective? (The behavior at the wall is hard to see, so you might need to reason it out)"
$ws3.Cells.Item(1,2).Value = "Correct"
$ws3.Cells.Item(1,3).Value = "Comment"
$ws3.Rows.Item(1).RowHeight = 135

$ws3.Cells.Item(2,1).Value = "Both gradients get steeper"
$ws3.Cells.Item(2,2).Value = "Y"
$ws3.Cells.Item(2,3).Value = "Yep!  The overall flux is increasing in both the solid and fluid.  That means that the temperature gradients in both must become steeper."
$ws3.Rows.Item(2).RowHeight = 75

$ws3.Cells.Item(3,1).Value = "The gradient in the solid gets steeper, but the one in the fluid becomes more shallow"
$ws3.Cells.Item(3,2).Value = "N"
$ws3.Cells.Item(3,3).Clear()
$ws3.Rows.Item(3).RowHeight = 60

$ws3.Cells.Item(4,1).Value = "Both gradients become more shallow"
$ws3.Cells.Item(4,2).Value = "N"
$ws3.Rows.Item(4).RowHeight = 30

$ws3.Cells.Item(5,1).Value = "The gradient in the solid becomes more shallow, bu the one in the fluid gets steeper"
$ws3.Cells.Item(5,2).Value = "N"
$ws3.Rows.Item(5).RowHeight = 60

# ---------------------------------------------------------------------
# Selection / active-sheet bookkeeping
# ---------------------------------------------------------------------
# "1_" loses tab-selection; its own selection collapses to A1:C5.
$ws1.Range("A1:C5").Select()

# "3_" becomes the active sheet/tab, with C3 selected.
$ws3.Activate()
$ws3.Range("C3").Select()
